$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45908
$ws.Range("B2").Value = 104.5
$ws.Range("C2").Value = 94.5
$ws.Range("D2").Value = 91.54000000000001
$ws.Range("E2").Value = 85.52
$ws.Range("F2").Value = 80.94
$ws.Range("G2").Value = 81.97
$ws.Range("H2").Value = 92.02
$ws.Range("I2").Value = 104.45
$ws.Range("J2").Value = 108.36
$ws.Range("K2").Value = 98.20999999999999
$ws.Range("L2").Value = 85.83
$ws.Range("M2").Value = 63
$ws.Range("N2").Value = 52.01
$ws.Range("O2").Value = 38.56
$ws.Range("P2").Value = 27.2
$ws.Range("Q2").Value = 27.08
$ws.Range("R2").Value = 27.2
$ws.Range("S2").Value = 48.64
$ws.Range("T2").Value = 73.08
$ws.Range("U2").Value = 91.15000000000001
$ws.Range("V2").Value = 98.20999999999999
$ws.Range("W2").Value = 95.26000000000001
$ws.Range("X2").Value = 90.51000000000001
$ws.Range("Y2").Value = 73.08
$ws.Range("Z2").Value = 76.37
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 94.02
$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 103.28
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 99.5
$ws.Range("AG2").Value = "11h-23h"
